# Updates the "cryptos" price list to the Fri Jun 23 08:46:03 UTC 2023
# GitHub Actions refresh: new Price (D) / Volume(1h) (E) figures for every
# coin row, plus Coin/Link/Price/Volume rewrites for the two row pairs that
# swapped rank order (41<->42: FraxShare/TheSandbox, 49<->50:
# EnergySwap/Decentraland).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a value as plain text, even when it looks numeric (e.g. "16.80" or
# "0.07200"), so trailing zeros are preserved exactly like the inline
# strings already on the sheet, instead of Excel normalising them to a
# number (16.8 / 0.072).
function Set-TextValue($cell, [string]$text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

$colIndex = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5 }

# row number => @{ column letter = new value }
$updates = [ordered]@{
    2 = @{ "D" = "29.986.53"; "E" = "  -0.31%  " }
    3 = @{ "D" = "1.884.43"; "E" = "  -1.16%  " }
    4 = @{ "D" = "0.9985"; "E" = "  -0.18%  " }
    5 = @{ "D" = "243.61"; "E" = "  -3.23%  " }
    6 = @{ "D" = "0.9987"; "E" = "  -0.13%  " }
    7 = @{ "D" = "0.4944"; "E" = "  -3.08%  " }
    8 = @{ "E" = "  -2.33%  " }
    9 = @{ "D" = "0.06653"; "E" = "  -2.14%  " }
    10 = @{ "D" = "1.879.88"; "E" = "  -1.36%  " }
    11 = @{ "D" = "16.80"; "E" = "  -2.74%  " }
    12 = @{ "D" = "0.07200"; "E" = "  -1.76%  " }
    13 = @{ "D" = "0.6702"; "E" = "  -3.93%  " }
    14 = @{ "D" = "86.56"; "E" = "  -0.01%  " }
    15 = @{ "D" = "4.892"; "E" = "  -0.40%  " }
    16 = @{ "D" = "29.946.41"; "E" = "  -0.41%  " }
    17 = @{ "D" = "0.000007869"; "E" = "  -3.88%  " }
    18 = @{ "D" = "0.9983"; "E" = "  -0.19%  " }
    19 = @{ "D" = "12.83"; "E" = "  -1.62%  " }
    20 = @{ "D" = "2.118.19"; "E" = "  -1.66%  " }
    21 = @{ "D" = "0.9977"; "E" = "  -0.20%  " }
    22 = @{ "D" = "4.793"; "E" = "  -0.71%  " }
    23 = @{ "D" = "5.895"; "E" = "  +2.78%  " }
    24 = @{ "D" = "9.127"; "E" = "  -1.54%  " }
    25 = @{ "D" = "150.41"; "E" = "  +1.87%  " }
    26 = @{ "D" = "142.79"; "E" = "  +5.80%  " }
    27 = @{ "D" = "17.07"; "E" = "  +0.05%  " }
    28 = @{ "D" = "1.929"; "E" = "  -3.29%  " }
    29 = @{ "D" = "1.388"; "E" = "  -1.30%  " }
    30 = @{ "D" = "4.225"; "E" = "  -0.78%  " }
    31 = @{ "D" = "0.08785"; "E" = "  -0.33%  " }
    32 = @{ "D" = "4.012"; "E" = "  +0.35%  " }
    33 = @{ "D" = "0.05060"; "E" = "  +0.16%  " }
    34 = @{ "D" = "0.7163"; "E" = "  -0.61%  " }
    35 = @{ "D" = "1.118"; "E" = "  -1.90%  " }
    36 = @{ "D" = "2.667"; "E" = "  -0.83%  " }
    37 = @{ "D" = "0.01801"; "E" = "  +6.35%  " }
    38 = @{ "D" = "2.698"; "E" = "  -4.05%  " }
    39 = @{ "D" = "2.180"; "E" = "  -3.88%  " }
    40 = @{ "D" = "0.9346"; "E" = "  -3.08%  " }
    41 = @{ "B" = "TheSandbox"; "C" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; "D" = "0.4244"; "E" = "  -1.54%  " }
    42 = @{ "B" = "FraxShare"; "C" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; "D" = "5.770"; "E" = "  -6.07%  " }
    43 = @{ "D" = "0.9989"; "E" = "  -0.01%  " }
    44 = @{ "D" = "103.28"; "E" = "  -1.36%  " }
    45 = @{ "D" = "7.428"; "E" = "  -2.43%  " }
    46 = @{ "D" = "0.1273"; "E" = "  -0.53%  " }
    47 = @{ "D" = "0.05678"; "E" = "  -1.07%  " }
    48 = @{ "D" = "32.68"; "E" = "  -2.23%  " }
    49 = @{ "B" = "Decentraland"; "C" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; "D" = "0.3778"; "E" = "  -0.99%  " }
    50 = @{ "B" = "EnergySwap"; "C" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; "D" = "8.273"; "E" = "  -2.08%  " }
    51 = @{ "D" = "56.25"; "E" = "  -1.23%  " }
}

foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    foreach ($col in $rowUpdates.Keys) {
        $cell = $ws.Cells.Item($row, $colIndex[$col])
        Set-TextValue $cell $rowUpdates[$col]
    }
}

Write-Host "Updated $($updates.Count) rows"
